$wb = $excel.ActiveWorkbook

# Update the "展览" (Exhibition) sheet
$wsExhibit = $wb.Worksheets.Item("展览")
$wsExhibit.Range("F2").Value = 1009
$wsExhibit.Range("F3").Value = 496

# Update the "全部类型" (All Types) sheet, which mirrors the exhibition data
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F2").Value = 1009
$wsAll.Range("F3").Value = 496
